# Patch the "Gun Load" event chain in the Events sheet:
#  - e050 ("Ammo Reload Order") gets reworded
#  - two brand-new events are inserted right after it: e050a ("No Gun Round
#    Loaded") and e050b ("Out of Main Gun Ammunition")
# NOTE: the order in which the B-column (text) cells are written matters -
# it controls the order new entries land in the shared-string table, which
# in turn is what the final saved workbook's <v> indices need to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: two fresh rows right after the existing e050 row (62).
$ws.Rows("63:64").Insert()

# Row 63 -> new event e050a.
$ws.Range("A63").Value = "e050a"
$bodyA = @"
<Bold>e050a No Gun Round Loaded</Bold> 
<InlineUIContainer><Button Content='r9.6' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
Since the gun is unloaded at the start of the round, choose one of the highlighted boxes to load the gun.
<LineBreak/><LineBreak/>
"@
$ws.Range("B63").Value = $bodyA

# Row 64 -> new event e050b (label only for now, body written further down).
$ws.Range("A64").Value = "e050b"

# Row 62 -> existing e050 event, reworded body text.
$bodyE050 = @"
<Bold>e050 Ammo Reload Order</Bold> 
<InlineUIContainer><Button Content='r4.73' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r5.23' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>  
<InlineUIContainer><Button Content='r9.6' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
Right click marker or one of blue boxes to to select from pull down menu where to set the Gun Reload marker and/or Ready Rack Ammo Reload marker. Alternative, right click on the Gun Load marker and choose from pull down menu.
<LineBreak/><LineBreak/>
"@
$ws.Range("B62").Value = $bodyE050

# Row 64 body text (written after B62 so shared-string ordering matches).
$bodyB = @"
<Bold>e050b Out of Main Gun Ammunition</Bold> 
<InlineUIContainer><Button Content='r9.6' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
Since the gun is unloaded at the start of the round, choose one of the highlighted boxes to load the gun.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Continue50b' Height='100' Width='100'></Image></InlineUIContainer>
"@
$ws.Range("B64").Value = $bodyB

# Row heights matching the new wrapped text.
$ws.Rows(62).RowHeight = 120
$ws.Rows(63).RowHeight = 75
$ws.Rows(64).RowHeight = 90

# Leave selection on the newly inserted e050a text cell.
$ws.Range("B63").Select()
